# The commit swaps the OOXML content of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml: the "Integral" theme (previously theme2.xml, the
# theme actually applied to the single Slide Master / the presentation
# design) becomes the default "Office Theme" colour palette, while the
# Office Theme content (previously theme1.xml, used only by the Notes
# Master) becomes "Integral".
#
# The PowerPoint object model exposes the colour palette of the theme
# that backs the Slide Master (and therefore the whole deck design)
# through Slide/SlideRange.ThemeColorScheme - a 12-slot scheme matching
# <a:clrScheme> (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink). We
# rewrite every slot to the standard Office theme RGB values, which is
# what a user does when switching the deck's Design from "Integral" to
# the built-in "Office Theme".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index -> (role, target "Office Theme" RGB)
# RGB() below packs 0xRRGGBB into PowerPoint's 0x00BBGGRR long.
function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$tcs.Colors(1).RGB  = RGB 0x00 0x00 0x00  # dk1
$tcs.Colors(2).RGB  = RGB 0xFF 0xFF 0xFF  # lt1
$tcs.Colors(3).RGB  = RGB 0x44 0x54 0x6A  # dk2
$tcs.Colors(4).RGB  = RGB 0xE7 0xE6 0xE6  # lt2
$tcs.Colors(5).RGB  = RGB 0x5B 0x9B 0xD5  # accent1
$tcs.Colors(6).RGB  = RGB 0xED 0x7D 0x31  # accent2
$tcs.Colors(7).RGB  = RGB 0xA5 0xA5 0xA5  # accent3
$tcs.Colors(8).RGB  = RGB 0xFF 0xC0 0x00  # accent4
$tcs.Colors(9).RGB  = RGB 0x44 0x72 0xC4  # accent5
$tcs.Colors(10).RGB = RGB 0x70 0xAD 0x47  # accent6
$tcs.Colors(11).RGB = RGB 0x05 0x63 0xC1  # hlink
$tcs.Colors(12).RGB = RGB 0x95 0x4F 0x72  # folHlink
